$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.346.08'
$ws.Range("E2").Value = '  +5.51%  '
$ws.Range("D3").Value = '1.921.72'
$ws.Range("E3").Value = '  +6.14%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '''253.96'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").Value = '''0.9970'
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("D7").Value = '''0.5178'
$ws.Range("E7").Value = '  +4.42%  '
$ws.Range("D8").Value = '''46.14'
$ws.Range("E8").Value = '  +7.08%  '
$ws.Range("D9").Value = '''0.2977'
$ws.Range("E9").Value = '  +6.14%  '
$ws.Range("D10").Value = '''0.06808'
$ws.Range("E10").Value = '  +6.45%  '
$ws.Range("D11").Value = '1.913.64'
$ws.Range("E11").Value = '  +5.73%  '
$ws.Range("D12").Value = '''17.62'
$ws.Range("E12").Value = '  +4.65%  '
$ws.Range("D13").Value = '''0.07312'
$ws.Range("E13").Value = '  +2.85%  '
$ws.Range("D14").Value = '''0.6899'
$ws.Range("E14").Value = '  +6.59%  '
$ws.Range("D15").Value = '''88.35'
$ws.Range("E15").Value = '  +7.48%  '
$ws.Range("D16").Value = '''4.939'
$ws.Range("E16").Value = '  +5.02%  '
$ws.Range("D17").Value = '30.345.41'
$ws.Range("E17").Value = '  +5.55%  '
$ws.Range("D18").Value = '''0.000007877'
$ws.Range("E18").Value = '  +6.93%  '
$ws.Range("D19").Value = '''0.9967'
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("E20").Value = '  +6.86%  '
$ws.Range("D21").Value = '2.162.71'
$ws.Range("E21").Value = '  +5.97%  '
$ws.Range("D22").Value = '''0.9963'
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("E23").Value = '  +5.70%  '
$ws.Range("D24").Value = '''5.765'
$ws.Range("E24").Value = '  +8.71%  '
$ws.Range("D25").Value = '''9.240'
$ws.Range("E25").Value = '  +4.05%  '
$ws.Range("D26").Value = '''140.08'
$ws.Range("E26").Value = '  +25.50%  '
$ws.Range("D27").Value = '''146.24'
$ws.Range("E27").Value = '  +2.31%  '
$ws.Range("D28").Value = '''17.41'
$ws.Range("E28").Value = '  +8.54%  '
$ws.Range("D29").Value = '''2.027'
$ws.Range("E29").Value = '  +7.65%  '
$ws.Range("D30").Value = '''1.385'
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("D31").Value = '''4.304'
$ws.Range("E31").Value = '  +2.85%  '
$ws.Range("D32").Value = '''0.08898'
$ws.Range("E32").Value = '  +6.37%  '
$ws.Range("D33").Value = '''4.052'
$ws.Range("E33").Value = '  +5.21%  '
$ws.Range("D34").Value = '''0.05148'
$ws.Range("E34").Value = '  +3.95%  '
$ws.Range("E35").Value = '  +6.20%  '
$ws.Range("D36").Value = '''0.7239'
$ws.Range("E36").Value = '  +7.75%  '
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("D38").Value = '''2.858'
$ws.Range("E38").Value = '  +8.36%  '
$ws.Range("D39").Value = '''2.331'
$ws.Range("E39").Value = '  +8.33%  '
$ws.Range("D40").Value = '''0.9732'
$ws.Range("E40").Value = '  +0.94%  '
$ws.Range("D41").Value = '''0.01701'
$ws.Range("E41").Value = '  +6.03%  '
$ws.Range("D42").Value = '''6.158'
$ws.Range("E42").Value = '  +3.22%  '
$ws.Range("D43").Value = '''0.4352'
$ws.Range("E43").Value = '  +5.41%  '
$ws.Range("E44").Value = '  +4.82%  '
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").Value = '''7.704'
$ws.Range("E46").Value = '  +6.47%  '
$ws.Range("D47").Value = '''0.1282'
$ws.Range("E47").Value = '  +4.68%  '
$ws.Range("D48").Value = '''0.05737'
$ws.Range("E48").Value = '  +4.49%  '
$ws.Range("D49").Value = '''8.597'
$ws.Range("E49").Value = '  +5.52%  '
$ws.Range("D50").Value = '''33.24'
$ws.Range("E50").Value = '  +6.11%  '
$ws.Range("D51").Value = '''0.3862'
$ws.Range("E51").Value = '  +6.71%  '
